$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 23.02.2022 22:00"

# Row 10 (EuroOil Opustena): change D10 from text "+0.3" to numeric 0.3
$ws.Range("D10").Value = 0.3

# Row 10: change E10 from text timestamp to a real Excel date/time value
# 2022-02-23 21:47:11 -> serial date 44615.9077662037
$ws.Range("E10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = (Get-Date -Year 2022 -Month 2 -Day 23 -Hour 21 -Minute 47 -Second 11)
